# Fixed a bug calcMulLevelRTP
# The data rows (2-20) on the active sheet were re-derived/re-ordered by the
# fixed calculation. Write the corrected values back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101,  9, 30, 15, 60, 15),
    @(901, 16, 15, 45, 60, 60),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(601,  9, 60, 67, 60, 42),
    @(801,  3, 67, 65, 52, 45),
    @(902,  1,  0,  0,  0,  0),
    @(1001,18, 30, 75, 60, 72),
    @(301,  6, 45, 30, 60, 45),
    @(501,  9, 52, 30, 75, 45),
    @(701,  3, 90, 45, 97, 15),
    @(201,  9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(401,  9, 48, 67, 75, 45),
    @(1101, 0, 15, 30, 30,  0),
    @(802,  0,  4,  5,  4,  0),
    @(2,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(502,  0,  4,  0,  0,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
